{"js": "// Resume edit: tighten the DocuSign bullet line from\n// \"Lead a small (5-person) Agile-development team responsible for ...\"\n// to \"Lead an Agile-development team responsible for ...\".\nconst oldText =\n  \"Lead a small (5-person) Agile-development team responsible for \" +\n  \"DocuSign's Developer Programs including creating DocuSign\\u2019s first \" +\n  \"Developer Center, SDKs, and improved API tools.\";\nconst newText =\n  \"Lead an Agile-development team responsible for \" +\n  \"DocuSign's Developer Programs including creating DocuSign\\u2019s first \" +\n  \"Developer Center, SDKs, and improved API tools.\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole (multi-run) span in one shot so Word collapses it\n  // back down to a single run carrying the shared formatting, matching\n  // how the author's edit merged the three original runs into one.\n  results.items[0].insertText(newText, \"Replace\");\n} else {\n  // Fallback: in case the text was already split/searched differently,\n  // try a narrower, still-unambiguous replacement of just the clause\n  // that actually changed.\n  const narrowOld = \"Lead a small (5-person) Agile-development team responsible for \";\n  const narrowNew = \"Lead an Agile-development team responsible for \";\n  const narrowResults = body.search(narrowOld, { matchCase: true });\n  await context.sync();\n  if (narrowResults.items.length > 0) {\n    narrowResults.items[0].insertText(narrowNew, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Resume edit: tighten the DocuSign bullet line from\n# \"Lead a small (5-person) Agile-development team responsible for ...\"\n# to \"Lead an Agile-development team responsible for ...\".\n\n$d = $word.ActiveDocument\n\n$rsquo = [char]0x2019\n\n$oldFull = \"Lead a small (5-person) Agile-development team responsible for \" +\n    \"DocuSign's Developer Programs including creating DocuSign\" + $rsquo +\n    \"s first Developer Center, SDKs, and improved API tools.\"\n$newFull = \"Lead an Agile-development team responsible for \" +\n    \"DocuSign's Developer Programs including creating DocuSign\" + $rsquo +\n    \"s first Developer Center, SDKs, and improved API tools.\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = $oldFull\n$found = $find.Execute($oldFull, $false, $false, $false, $false, $false, $true, 1, $false, $newFull, 2)\n\nif (-not $found) {\n    # Fallback: narrower replace of just the clause that changed, in case\n    # the long phrase above doesn't match verbatim for some reason.\n    $narrowOld = \"Lead a small (5-person) Agile-development team responsible for \"\n    $narrowNew = \"Lead an Agile-development team responsible for \"\n    $rng2 = $d.Content\n    $find2 = $rng2.Find\n    $find2.Text = $narrowOld\n    $find2.Execute($narrowOld, $false, $false, $false, $false, $false, $true, 1, $false, $narrowNew, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
